$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(79, 1).Value = 46028
$ws.Cells.Item(79, 2).Value = 176
$ws.Cells.Item(79, 3).Value = 184
$ws.Cells.Item(79, 4).Value = 181

$ws.Cells.Item(79, 1).NumberFormat = $ws.Cells.Item(78, 1).NumberFormat
